$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has headers Nama | Waktu | Jenis in columns A-C, with Arena
# rows below. We add an "Audio" column, pushing the existing "Jenis" column (and
# its Tanding/Tunggal/Ganda/Regu values) one column to the right, into column D.

# Remember the current "Jenis" column's values (row 1 header + the 4 data rows)
# before they get overwritten.
$jenisHeader = $ws.Cells.Item(1, 3).Text
$jenisRow2 = $ws.Cells.Item(2, 3).Text
$jenisRow3 = $ws.Cells.Item(3, 3).Text
$jenisRow4 = $ws.Cells.Item(4, 3).Text
$jenisRow5 = $ws.Cells.Item(5, 3).Text

# Column D should look exactly like column C did (same header styling), so copy
# the formatting over before writing the relocated values into it.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(1, 4).Value = $jenisHeader
$ws.Cells.Item(2, 4).Value = $jenisRow2
$ws.Cells.Item(3, 4).Value = $jenisRow3
$ws.Cells.Item(4, 4).Value = $jenisRow4
$ws.Cells.Item(5, 4).Value = $jenisRow5

# Column C becomes the new "Audio" column, referencing the audio clip file name.
$ws.Cells.Item(1, 3).Value = "Audio"
$ws.Cells.Item(2, 3).Value = "audio.mp3"
$ws.Cells.Item(3, 3).Value = "audio.mp3"
$ws.Cells.Item(4, 3).Value = "audio.mp3"
$ws.Cells.Item(5, 3).Value = "audio.mp3"
